# Update localization status report: two files (27f1d26c..., caf09c3d...)
# moved from "Ready for handoff" to "In Translation" state.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "In Translation"   # zh-cn column, 27f1d26c...md row
$wsOverview.Range("C3").Value = "In Translation"   # de-de column, 27f1d26c...md row
$wsOverview.Range("B4").Value = "In Translation"   # zh-cn column, caf09c3d...md row
$wsOverview.Range("C4").Value = "In Translation"   # de-de column, caf09c3d...md row

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "In Translation"       # Status, 27f1d26c...md row
$wsZhCn.Range("C4").Value = "In Translation"       # Status, caf09c3d...md row

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "In Translation"       # Status, 27f1d26c...md row
$wsDeDe.Range("C4").Value = "In Translation"       # Status, caf09c3d...md row
